# Add new simulation results (columns V:AO) to Sheet1, mirroring the
# existing angle row (row 2) and appending a freshly generated 0/1 row (row 3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (angles, columns V through AO) - repeats the same 0..2*pi sweep
# that already occupies columns B through U.
$row2Values = @(0.31415926535897898, 0.62831853071795896, 0.94247779607693805, 1.2566370614359199, 1.5707963267949001, 1.8849555921538801, 2.1991148575128601, 2.5132741228718301, 2.8274333882308098, 3.14159265358979, 3.4557519189487702, 3.76991118430775, 4.0840704496667302, 4.3982297150257104, 4.7123889803846897, 5.0265482457436699, 5.3407075111026501, 5.6548667764616303, 5.9690260418206096, 6.2831853071795898)

# Row 3 (new simulation outcomes, columns V through AO)
$row3Values = @(1, 0, 0, 0, 0, 1, 1, 0, 1, 0, 1, 1, 0, 1, 0, 1, 0, 0, 1, 0)

$startCol = 22  # column V

for ($i = 0; $i -lt $row2Values.Length; $i++) {
    $col = $startCol + $i
    $ws.Cells.Item(2, $col).Value = $row2Values[$i]
    $ws.Cells.Item(3, $col).Value = $row3Values[$i]
}

# Update the view: scroll so column L is at the left edge and select AA5,
# matching the saved selection/scroll state from the authoring session.
$ws.Range("L1").Select()
$excel.ActiveWindow.ScrollColumn = 12
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AA5").Select()
